$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 34

# Column A holds a date-like label ("01-07-2021"). Excel's COM layer will
# auto-convert a bare string like this into a real date serial + date
# number format, so force the cell to text first (matching the sibling
# cells above it, which are plain shared-string text, not dates).
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "01-07-2021"
# Restore the same (default) style as the rest of the column so no stray
# per-cell style index is left behind.
$ws.Cells.Item($row, 1).Style = $ws.Cells.Item($row - 1, 1).Style

$ws.Cells.Item($row, 2).Value = 202
$ws.Cells.Item($row, 3).Value = 50
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 152
$ws.Cells.Item($row, 11).Value = 0
